# The workbook contains a weekly price log for "Brócoli" at the
# "Terminal Hortofrutícola Agro Chillán" market, sorted (roughly) by date.
# A week that was missing from the log (2021-09-21, serial 44460) is being
# added back in. It belongs right after the existing 2021-08-10 entry
# (row 101) / before the 2021-08-06 entry (row 102), so every row from the
# old row 102 down to the old row 160 shifts down by one, and the new row
# is populated with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 102:160 down to 103:161, duplicating formatting (incl. the
# date number-format on column D) from the row being split.
$ws.Rows.Item(102).Insert()

# Populate the newly-inserted row with the missing week's record.
$ws.Range("A102").Value = 7
$ws.Range("B102").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value = "Ñuble"
$ws.Range("D102").Value = 44460
$ws.Range("E102").Value = 16
$ws.Range("F102").Value = 100112023
$ws.Range("G102").Value = "Brócoli"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 160
$ws.Range("K102").Value = 700
$ws.Range("L102").Value = 750
$ws.Range("M102").Value = 725
$ws.Range("N102").Value = "$/unidad"
$ws.Range("O102").Value = "Región del Maule"
$ws.Range("P102").Value = 725
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"
